$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 790
$wsExhibit.Range("F3").Value = 62
$wsExhibit.Range("F5").Value = 148
$wsExhibit.Range("G5").Value = 169
$wsExhibit.Range("F7").Value = 169
$wsExhibit.Range("F9").Value = 467
$wsExhibit.Range("F10").Value = 519
$wsExhibit.Range("F11").Value = 147
$wsExhibit.Range("F12").Value = 11975
$wsExhibit.Range("F13").Value = 5434

# Sheet "全部类型" (All Types) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 790
$wsAll.Range("F3").Value = 62
$wsAll.Range("F7").Value = 148
$wsAll.Range("G7").Value = 169
$wsAll.Range("F9").Value = 169
$wsAll.Range("F11").Value = 467
$wsAll.Range("F12").Value = 519
$wsAll.Range("F13").Value = 147
$wsAll.Range("F14").Value = 11975
$wsAll.Range("F16").Value = 5434
